# This script reproduces a pandas re-export where `df.to_excel(...)` was
# changed from an implicit/explicit index (year as the DataFrame index,
# written via `set_index`) to `index=False`. The visible effect in the
# workbook is:
#   1. On the main "df_all" sheet, the leading "year" column (A) is removed
#      entirely and every other column shifts one place to the left
#      (so the old M column - "...65 years and over" - becomes the new L).
#   2. On every per-category sheet ("Total fatalities", "Private industrys",
#      ... "65 years and over"), the "year" column is still present (those
#      sheets keep their own small 2-column frame), but the bold/bordered
#      header style that used to bleed down into the data rows of column A
#      is cleared from the data rows, leaving only the header cell (row 1)
#      styled.

$wb = $excel.ActiveWorkbook

# 1) df_all: drop the leading "year" index column so the data shifts left.
$dfAll = $wb.Worksheets.Item("df_all")
$dfAll.Columns("A").Delete()

# 2) Per-category sheets: clear the inherited header formatting from the
#    "year" data cells (row 2 downward), keeping the header cell styled.
$sheetNames = @(
    "Total fatalities",
    "Private industrys",
    "General warehousing and storage",
    "Under 16 years",
    "16 to 17 years",
    "18 to 19 years",
    "20 to 24 years",
    "25 to 34 years",
    "35 to 44 years",
    "45 to 54 years",
    "55 to 64 years",
    "65 years and over"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $lastRow = $ws.UsedRange.Rows.Count
    if ($lastRow -ge 2) {
        $ws.Range("A2:A" + $lastRow).ClearFormats()
    }
}
